{"js": "// Bullet Patterns: Decided on standard Base/Child inheritance.\n// \"Underscore Prefix\" -> two runs: the quoted \"m_\" literal, then \"Prefix\".\nconst body = context.document.body;\n\nconst results = body.search(\"Underscore Prefix\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n\n  // Insert the two new runs (identical Times New Roman formatting, but as\n  // separate <w:r> elements) right before the existing text, while that\n  // text is still present so the paragraph's properties stay untouched.\n  const insertionPoint = target.getRange(\"Start\");\n\n  const quote = \"\\u201C\";\n  const closeQuote = \"\\u201D\";\n  const xml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' +\n    '<w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/></w:rPr>' +\n    '<w:t xml:space=\"preserve\">' + quote + \"m_\" + closeQuote + ' </w:t></w:r>' +\n    '<w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/></w:rPr>' +\n    \"<w:t>Prefix</w:t></w:r>\" +\n    \"</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\n  insertionPoint.insertOoxml(xml, \"Start\");\n  await context.sync();\n\n  // Remove the original \"Underscore Prefix\" text. Using delete() (rather\n  // than overwriting text in place) avoids the host re-merging the freshly\n  // inserted runs with whatever used to sit next to them.\n  const stale = body.search(\"Underscore Prefix\", { matchCase: true });\n  stale.load(\"items\");\n  await context.sync();\n\n  if (stale.items.length > 0) {\n    stale.items[0].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Bullet Patterns: Decided on standard Base/Child inheritance.\n# \"Underscore Prefix\" -> two runs: the quoted \"m_\" literal, then \"Prefix\".\n$d = $word.ActiveDocument\n\n# Locate the paragraph's text to edit.\n$find = $d.Content.Find\n$find.Text = \"Underscore Prefix\"\n$find.Execute()\n\nif ($find.Found) {\n    $target = $find.Parent\n\n    # Insert the two new runs (identical Times New Roman formatting, split\n    # into separate <w:r> elements) immediately before the existing text,\n    # while the paragraph still has content so its <w:pPr> is untouched.\n    $insertionPoint = $d.Content\n    $insertionPoint.Start = $target.Start\n    $insertionPoint.End = $target.Start\n\n    $quote = [char]0x201C\n    $closeQuote = [char]0x201D\n    $xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/></w:rPr><w:t xml:space=\"preserve\">' + $quote + 'm_' + $closeQuote + ' </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/></w:rPr><w:t>Prefix</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $insertionPoint.InsertXML($xml)\n\n    # Remove the original \"Underscore Prefix\" run text. Using Delete()\n    # (rather than assigning Text = \"\") avoids the host re-merging the\n    # freshly inserted runs with whatever used to sit next to them.\n    $find2 = $d.Content.Find\n    $find2.Text = \"Underscore Prefix\"\n    $find2.Execute()\n    if ($find2.Found) {\n        $find2.Parent.Delete()\n    }\n}\n"}
